$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 5 new blank rows above row 5. This pushes the former rows 5-12
#    down to become rows 10-17 (Excel keeps their values/styles/formulas).
# ---------------------------------------------------------------------------
$ws.Rows("5:9").Insert()

# ---------------------------------------------------------------------------
# 2. Populate the 5 new rows with the "COMMON" skill (a non-skill placeholder
#    used by the new SkillAttributeManager to hold attributes common to every
#    skill: ATTACK, CRITICAL_RATE, CRITICAL_AMOUNT, AREA, ATTACK_DECREASE).
# ---------------------------------------------------------------------------

# Row 5 - skill header columns (id / name / description / skillType)
$c = $ws.Range("B5")
$c.Style = "Normal"
$c.VerticalAlignment = -4108
$c.Value = "COMMON"

$c = $ws.Range("C5")
$c.Style = "Normal"
$c.Value = "COMMON"

$c = $ws.Range("D5")
$c.Style = "Normal"
$c.Value = "Not a skill"

$c = $ws.Range("E5")
$c.Style = "Normal"
$c.Value = "PASSIVE"

$ws.Range("F5").Style = "Normal"

# Row 5 - first attribute row (ATTACK / 1)
$c = $ws.Range("G5")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = "ATTACK"

$c = $ws.Range("H5")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = 1

$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Style = "Normal"
$ws.Range("L5").Style = "Normal"

# Row 6 - CRITICAL_RATE / 0.05
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("F6").Style = "Normal"

$c = $ws.Range("G6")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = "CRITICAL_RATE"

$c = $ws.Range("H6")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = 0.05

$ws.Range("I6").Style = "Normal"
$ws.Range("J6").Style = "Normal"
$ws.Range("L6").Style = "Normal"

# Row 7 - CRITICAL_AMOUNT / 2
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("F7").Style = "Normal"

$c = $ws.Range("G7")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = "CRITICAL_AMOUNT"

$c = $ws.Range("H7")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = 2

$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Style = "Normal"
$ws.Range("L7").Style = "Normal"

# Row 8 - AREA / 1
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("F8").Style = "Normal"

$c = $ws.Range("G8")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = "AREA"

$c = $ws.Range("H8")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = 1

$ws.Range("I8").Style = "Normal"
$ws.Range("J8").Style = "Normal"
$ws.Range("L8").Style = "Normal"

# Row 9 - ATTACK_DECREASE / 1
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("F9").Style = "Normal"

$c = $ws.Range("G9")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = "ATTACK_DECREASE"

$c = $ws.Range("H9")
$c.Style = "Normal"
$c.HorizontalAlignment = -4108
$c.Value = 1

$ws.Range("I9").Style = "Normal"
$ws.Range("J9").Style = "Normal"
$ws.Range("L9").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. The Chain Lightning skill (now on row 10) no longer depends on
#    EXPLOSION;SLOW - clear the dependencies cell.
# ---------------------------------------------------------------------------
$ws.Range("F10").ClearContents()

# ---------------------------------------------------------------------------
# 4. Restore the view/selection like the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("F14").Select()
